$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: mooring serial number + deployment number ---
$ws1.Range("A2").Value = "CP05MOAS-GL376"
$ws1.Range("C2").Value = 1

# --- Asset_Cal_Info sheet: instrument reference designators renamed GL003 -> GL376 ---
# and deployment number 2 -> 1 for every populated row.

# ADCPAM000 block (rows 2-5)
$ws2.Range("A2").Value = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("C2").Value = 1
$ws2.Range("A3").Value = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("C3").Value = 1
$ws2.Range("A4").Value = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("C4").Value = 1
$ws2.Range("A5").Value = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("C5").Value = 1

# FLORTM000 block (rows 7-10)
$ws2.Range("A7").Value = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("C7").Value = 1
$ws2.Range("A8").Value = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("C8").Value = 1
$ws2.Range("A9").Value = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("C9").Value = 1
$ws2.Range("A10").Value = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("C10").Value = 1

# CTDGVM000 (row 12)
$ws2.Range("A12").Value = "CP05MOAS-GL376-03-CTDGVM000"
$ws2.Range("C12").Value = 1

# DOSTAM000 (row 14)
$ws2.Range("A14").Value = "CP05MOAS-GL376-04-DOSTAM000"
$ws2.Range("C14").Value = 1

# PARADM000 (row 16)
$ws2.Range("A16").Value = "CP05MOAS-GL376-05-PARADM000"
$ws2.Range("C16").Value = 1

# ENG000000 (row 18)
$ws2.Range("A18").Value = "CP05MOAS-GL376-00-ENG000000"
$ws2.Range("C18").Value = 1

# --- View / selection state ---
# Moorings tab loses the "active" flag; its lingering selection moves to D19.
$ws1.Range("D19").Select()

# Asset_Cal_Info becomes the active/visible tab, selection moves to C20.
$ws2.Activate()
$ws2.Range("C20").Select()
